$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.037445187568665
$ws.Range("B1").Value = 2.173807382583618
$ws.Range("C1").Value = 4.051377296447754
$ws.Range("D1").Value = 0.830157458782196
$ws.Range("E1").Value = 0.9154462814331055
